$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = [double]"16.60580966666667"
$ws.Range("H2").Value = [double]"49.817429"
$ws.Range("I2").Value = [double]"0.8577430056934899"
$ws.Range("J2").Value = [double]"0.8577430056934898"
$ws.Range("M2").Value = [double]"11.195312"
$ws.Range("N2").Value = [double]"33.585936"
$ws.Range("O2").Value = [double]"0.8263729865424168"
$ws.Range("P2").Value = [double]"0.8263729865424166"
$ws.Range("Q2").Value = [double]"185.9072202309494"
$ws.Range("R2").Value = [double]"1673.164982078544"
$ws.Range("S2").Value = [double]"0.7088156493007984"
$ws.Range("T2").Value = [double]"0.7088156493007983"
$ws.Range("G3").Value = [double]"16.60580966666667"
$ws.Range("H3").Value = [double]"49.817429"
$ws.Range("I3").Value = [double]"0.8577430056934899"
$ws.Range("J3").Value = [double]"0.8577430056934898"
$ws.Range("O3").Value = [double]"0.1468201815007494"
$ws.Range("P3").Value = [double]"0.1468201815007494"
$ws.Range("Q3").Value = [double]"33.02979678802312"
$ws.Range("R3").Value = [double]"297.268171092208"
$ws.Range("S3").Value = [double]"0.1259339837769166"
$ws.Range("T3").Value = [double]"0.1259339837769165"
$ws.Range("G4").Value = [double]"16.60580966666667"
$ws.Range("H4").Value = [double]"49.817429"
$ws.Range("I4").Value = [double]"0.8577430056934899"
$ws.Range("J4").Value = [double]"0.8577430056934898"
$ws.Range("K4").Value = [double]"2"
$ws.Range("L4").Value = [double]"0.6666666666666666"
$ws.Range("M4").Value = [double]"0.02602966666666666"
$ws.Range("N4").Value = [double]"0.07808899999999999"
$ws.Range("O4").Value = [double]"0.001921358992231473"
$ws.Range("P4").Value = [double]"0.001921358992231473"
$ws.Range("Q4").Value = [double]"0.4322436903534445"
$ws.Range("R4").Value = [double]"3.890193213181"
$ws.Range("S4").Value = [double]"0.001648032237012839"
$ws.Range("T4").Value = [double]"0.001648032237012839"
$ws.Range("G5").Value = [double]"16.60580966666667"
$ws.Range("H5").Value = [double]"49.817429"
$ws.Range("I5").Value = [double]"0.8577430056934899"
$ws.Range("J5").Value = [double]"0.8577430056934898"
$ws.Range("M5").Value = [double]"0.2855723333333333"
$ws.Range("N5").Value = [double]"0.856717"
$ws.Range("O5").Value = [double]"0.02107929300858727"
$ws.Range("P5").Value = [double]"0.02107929300858727"
$ws.Range("Q5").Value = [double]"4.742159813399223"
$ws.Range("R5").Value = [double]"42.679438320593"
$ws.Range("S5").Value = [double]"0.01808061614307941"
$ws.Range("T5").Value = [double]"0.01808061614307941"
$ws.Range("G6").Value = [double]"16.60580966666667"
$ws.Range("H6").Value = [double]"49.817429"
$ws.Range("I6").Value = [double]"0.8577430056934899"
$ws.Range("J6").Value = [double]"0.8577430056934898"
$ws.Range("K6").Value = [double]"2"
$ws.Range("L6").Value = [double]"0.6666666666666666"
$ws.Range("M6").Value = [double]"0.04091133333333333"
$ws.Range("N6").Value = [double]"0.122734"
$ws.Range("O6").Value = [double]"0.003019837295298156"
$ws.Range("P6").Value = [double]"0.003019837295298156"
$ws.Range("Q6").Value = [double]"0.679365814542889"
$ws.Range("R6").Value = [double]"6.114292330886001"
$ws.Range("S6").Value = [double]"0.00259024431837434"
$ws.Range("T6").Value = [double]"0.002590244318374339"
$ws.Range("G7").Value = [double]"16.60580966666667"
$ws.Range("H7").Value = [double]"49.817429"
$ws.Range("I7").Value = [double]"0.8577430056934899"
$ws.Range("J7").Value = [double]"0.8577430056934898"
$ws.Range("K7").Value = [double]"1"
$ws.Range("L7").Value = [double]"0.3333333333333333"
$ws.Range("M7").Value = [double]"0.010653"
$ws.Range("N7").Value = [double]"0.031959"
$ws.Range("O7").Value = [double]"0.000786342660716947"
$ws.Range("P7").Value = [double]"0.0007863426607169469"
$ws.Range("Q7").Value = [double]"0.1769016903790001"
$ws.Range("R7").Value = [double]"1.592115213411"
$ws.Range("S7").Value = [double]"0.0006744799173083703"
$ws.Range("T7").Value = [double]"0.0006744799173083701"
$ws.Range("I8").Value = [double]"0.1296644119935396"
$ws.Range("J8").Value = [double]"0.1296644119935396"
$ws.Range("M8").Value = [double]"11.195312"
$ws.Range("N8").Value = [double]"33.585936"
$ws.Range("O8").Value = [double]"0.8263729865424168"
$ws.Range("P8").Value = [double]"0.8263729865424166"
$ws.Range("Q8").Value = [double]"28.10346483339733"
$ws.Range("R8").Value = [double]"252.931183500576"
$ws.Range("S8").Value = [double]"0.1071511673873677"
$ws.Range("T8").Value = [double]"0.1071511673873677"
$ws.Range("I9").Value = [double]"0.1296644119935396"
$ws.Range("J9").Value = [double]"0.1296644119935396"
$ws.Range("O9").Value = [double]"0.1468201815007494"
$ws.Range("P9").Value = [double]"0.1468201815007494"
$ws.Range("S9").Value = [double]"0.01903735250307944"
$ws.Range("T9").Value = [double]"0.01903735250307944"
$ws.Range("I10").Value = [double]"0.1296644119935396"
$ws.Range("J10").Value = [double]"0.1296644119935396"
$ws.Range("K10").Value = [double]"2"
$ws.Range("L10").Value = [double]"0.6666666666666666"
$ws.Range("M10").Value = [double]"0.02602966666666666"
$ws.Range("N10").Value = [double]"0.07808899999999999"
$ws.Range("O10").Value = [double]"0.001921358992231473"
$ws.Range("P10").Value = [double]"0.001921358992231473"
$ws.Range("Q10").Value = [double]"0.06534197723044442"
$ws.Range("R10").Value = [double]"0.5880777950739999"
$ws.Range("S10").Value = [double]"0.0002491318839561939"
$ws.Range("T10").Value = [double]"0.0002491318839561939"
$ws.Range("I11").Value = [double]"0.1296644119935396"
$ws.Range("J11").Value = [double]"0.1296644119935396"
$ws.Range("M11").Value = [double]"0.2855723333333333"
$ws.Range("N11").Value = [double]"0.856717"
$ws.Range("O11").Value = [double]"0.02107929300858727"
$ws.Range("P11").Value = [double]"0.02107929300858727"
$ws.Range("Q11").Value = [double]"0.7168689918802221"
$ws.Range("R11").Value = [double]"6.451820926921999"
$ws.Range("S11").Value = [double]"0.002733234133197999"
$ws.Range("T11").Value = [double]"0.002733234133197999"
$ws.Range("I12").Value = [double]"0.1296644119935396"
$ws.Range("J12").Value = [double]"0.1296644119935396"
$ws.Range("K12").Value = [double]"2"
$ws.Range("L12").Value = [double]"0.6666666666666666"
$ws.Range("M12").Value = [double]"0.04091133333333333"
$ws.Range("N12").Value = [double]"0.122734"
$ws.Range("O12").Value = [double]"0.003019837295298156"
$ws.Range("P12").Value = [double]"0.003019837295298156"
$ws.Range("Q12").Value = [double]"0.1026992564048889"
$ws.Range("R12").Value = [double]"0.924293307644"
$ws.Range("S12").Value = [double]"0.0003915654272109965"
$ws.Range("T12").Value = [double]"0.0003915654272109964"
$ws.Range("I13").Value = [double]"0.1296644119935396"
$ws.Range("J13").Value = [double]"0.1296644119935396"
$ws.Range("K13").Value = [double]"1"
$ws.Range("L13").Value = [double]"0.3333333333333333"
$ws.Range("M13").Value = [double]"0.010653"
$ws.Range("N13").Value = [double]"0.031959"
$ws.Range("O13").Value = [double]"0.000786342660716947"
$ws.Range("P13").Value = [double]"0.0007863426607169469"
$ws.Range("Q13").Value = [double]"0.026742105166"
$ws.Range("R13").Value = [double]"0.240678946494"
$ws.Range("S13").Value = [double]"0.0001019606587272984"
$ws.Range("T13").Value = [double]"0.0001019606587272983"
$ws.Range("G14").Value = [double]"0.21333"
$ws.Range("H14").Value = [double]"0.6399899999999999"
$ws.Range("I14").Value = [double]"0.01101917455864245"
$ws.Range("J14").Value = [double]"0.01101917455864245"
$ws.Range("M14").Value = [double]"11.195312"
$ws.Range("N14").Value = [double]"33.585936"
$ws.Range("O14").Value = [double]"0.8263729865424168"
$ws.Range("P14").Value = [double]"0.8263729865424166"
$ws.Range("Q14").Value = [double]"2.38829590896"
$ws.Range("R14").Value = [double]"21.49466318064"
$ws.Range("S14").Value = [double]"0.009105948189257577"
$ws.Range("T14").Value = [double]"0.009105948189257576"
$ws.Range("G15").Value = [double]"0.21333"
$ws.Range("H15").Value = [double]"0.6399899999999999"
$ws.Range("I15").Value = [double]"0.01101917455864245"
$ws.Range("J15").Value = [double]"0.01101917455864245"
$ws.Range("O15").Value = [double]"0.1468201815007494"
$ws.Range("P15").Value = [double]"0.1468201815007494"
$ws.Range("Q15").Value = [double]"0.42432417872"
$ws.Range("R15").Value = [double]"3.81891760848"
$ws.Range("S15").Value = [double]"0.001617837208688325"
$ws.Range("T15").Value = [double]"0.001617837208688325"
$ws.Range("G16").Value = [double]"0.21333"
$ws.Range("H16").Value = [double]"0.6399899999999999"
$ws.Range("I16").Value = [double]"0.01101917455864245"
$ws.Range("J16").Value = [double]"0.01101917455864245"
$ws.Range("K16").Value = [double]"2"
$ws.Range("L16").Value = [double]"0.6666666666666666"
$ws.Range("M16").Value = [double]"0.02602966666666666"
$ws.Range("N16").Value = [double]"0.07808899999999999"
$ws.Range("O16").Value = [double]"0.001921358992231473"
$ws.Range("P16").Value = [double]"0.001921358992231473"
$ws.Range("Q16").Value = [double]"0.005552908789999999"
$ws.Range("R16").Value = [double]"0.04997617910999999"
$ws.Range("S16").Value = [double]"2.117179012521594E-05"
$ws.Range("T16").Value = [double]"2.117179012521594E-05"
$ws.Range("G17").Value = [double]"0.21333"
$ws.Range("H17").Value = [double]"0.6399899999999999"
$ws.Range("I17").Value = [double]"0.01101917455864245"
$ws.Range("J17").Value = [double]"0.01101917455864245"
$ws.Range("M17").Value = [double]"0.2855723333333333"
$ws.Range("N17").Value = [double]"0.856717"
$ws.Range("O17").Value = [double]"0.02107929300858727"
$ws.Range("P17").Value = [double]"0.02107929300858727"
$ws.Range("Q17").Value = [double]"0.06092114586999999"
$ws.Range("R17").Value = [double]"0.54829031283"
$ws.Range("S17").Value = [double]"0.0002322764092343945"
$ws.Range("T17").Value = [double]"0.0002322764092343945"
$ws.Range("G18").Value = [double]"0.21333"
$ws.Range("H18").Value = [double]"0.6399899999999999"
$ws.Range("I18").Value = [double]"0.01101917455864245"
$ws.Range("J18").Value = [double]"0.01101917455864245"
$ws.Range("K18").Value = [double]"2"
$ws.Range("L18").Value = [double]"0.6666666666666666"
$ws.Range("M18").Value = [double]"0.04091133333333333"
$ws.Range("N18").Value = [double]"0.122734"
$ws.Range("O18").Value = [double]"0.003019837295298156"
$ws.Range("P18").Value = [double]"0.003019837295298156"
$ws.Range("Q18").Value = [double]"0.00872761474"
$ws.Range("R18").Value = [double]"0.07854853266"
$ws.Range("S18").Value = [double]"3.327611429558906E-05"
$ws.Range("T18").Value = [double]"3.327611429558906E-05"
$ws.Range("G19").Value = [double]"0.21333"
$ws.Range("H19").Value = [double]"0.6399899999999999"
$ws.Range("I19").Value = [double]"0.01101917455864245"
$ws.Range("J19").Value = [double]"0.01101917455864245"
$ws.Range("K19").Value = [double]"1"
$ws.Range("L19").Value = [double]"0.3333333333333333"
$ws.Range("M19").Value = [double]"0.010653"
$ws.Range("N19").Value = [double]"0.031959"
$ws.Range("O19").Value = [double]"0.000786342660716947"
$ws.Range("P19").Value = [double]"0.0007863426607169469"
$ws.Range("Q19").Value = [double]"0.00227260449"
$ws.Range("R19").Value = [double]"0.02045344041"
$ws.Range("S19").Value = [double]"8.664847041347392E-06"
$ws.Range("T19").Value = [double]"8.664847041347391E-06"
$ws.Range("G20").Value = [double]"0.030461"
$ws.Range("H20").Value = [double]"0.09138300000000001"
$ws.Range("I20").Value = [double]"0.001573407754328072"
$ws.Range("J20").Value = [double]"0.001573407754328072"
$ws.Range("M20").Value = [double]"11.195312"
$ws.Range("N20").Value = [double]"33.585936"
$ws.Range("O20").Value = [double]"0.8263729865424168"
$ws.Range("P20").Value = [double]"0.8263729865424166"
$ws.Range("Q20").Value = [double]"0.3410203988320001"
$ws.Range("R20").Value = [double]"3.069183589488"
$ws.Range("S20").Value = [double]"0.001300221664993086"
$ws.Range("T20").Value = [double]"0.001300221664993086"
$ws.Range("G21").Value = [double]"0.030461"
$ws.Range("H21").Value = [double]"0.09138300000000001"
$ws.Range("I21").Value = [double]"0.001573407754328072"
$ws.Range("J21").Value = [double]"0.001573407754328072"
$ws.Range("O21").Value = [double]"0.1468201815007494"
$ws.Range("P21").Value = [double]"0.1468201815007494"
$ws.Range("Q21").Value = [double]"0.06058847235733334"
$ws.Range("R21").Value = [double]"0.5452962512160001"
$ws.Range("S21").Value = [double]"0.0002310080120651342"
$ws.Range("T21").Value = [double]"0.0002310080120651342"
$ws.Range("G22").Value = [double]"0.030461"
$ws.Range("H22").Value = [double]"0.09138300000000001"
$ws.Range("I22").Value = [double]"0.001573407754328072"
$ws.Range("J22").Value = [double]"0.001573407754328072"
$ws.Range("K22").Value = [double]"2"
$ws.Range("L22").Value = [double]"0.6666666666666666"
$ws.Range("M22").Value = [double]"0.02602966666666666"
$ws.Range("N22").Value = [double]"0.07808899999999999"
$ws.Range("O22").Value = [double]"0.001921358992231473"
$ws.Range("P22").Value = [double]"0.001921358992231473"
$ws.Range("Q22").Value = [double]"0.0007928896763333333"
$ws.Range("R22").Value = [double]"0.007136007087"
$ws.Range("S22").Value = [double]"3.023081137224971E-06"
$ws.Range("T22").Value = [double]"3.023081137224971E-06"
$ws.Range("G23").Value = [double]"0.030461"
$ws.Range("H23").Value = [double]"0.09138300000000001"
$ws.Range("I23").Value = [double]"0.001573407754328072"
$ws.Range("J23").Value = [double]"0.001573407754328072"
$ws.Range("M23").Value = [double]"0.2855723333333333"
$ws.Range("N23").Value = [double]"0.856717"
$ws.Range("O23").Value = [double]"0.02107929300858727"
$ws.Range("P23").Value = [double]"0.02107929300858727"
$ws.Range("Q23").Value = [double]"0.008698818845666667"
$ws.Range("R23").Value = [double]"0.078289369611"
$ws.Range("S23").Value = [double]"3.316632307546473E-05"
$ws.Range("T23").Value = [double]"3.316632307546472E-05"
$ws.Range("G24").Value = [double]"0.030461"
$ws.Range("H24").Value = [double]"0.09138300000000001"
$ws.Range("I24").Value = [double]"0.001573407754328072"
$ws.Range("J24").Value = [double]"0.001573407754328072"
$ws.Range("K24").Value = [double]"2"
$ws.Range("L24").Value = [double]"0.6666666666666666"
$ws.Range("M24").Value = [double]"0.04091133333333333"
$ws.Range("N24").Value = [double]"0.122734"
$ws.Range("O24").Value = [double]"0.003019837295298156"
$ws.Range("P24").Value = [double]"0.003019837295298156"
$ws.Range("Q24").Value = [double]"0.001246200124666667"
$ws.Range("R24").Value = [double]"0.011215801122"
$ws.Range("S24").Value = [double]"4.751435417231232E-06"
$ws.Range("T24").Value = [double]"4.751435417231231E-06"
$ws.Range("G25").Value = [double]"0.030461"
$ws.Range("H25").Value = [double]"0.09138300000000001"
$ws.Range("I25").Value = [double]"0.001573407754328072"
$ws.Range("J25").Value = [double]"0.001573407754328072"
$ws.Range("K25").Value = [double]"1"
$ws.Range("L25").Value = [double]"0.3333333333333333"
$ws.Range("M25").Value = [double]"0.010653"
$ws.Range("N25").Value = [double]"0.031959"
$ws.Range("O25").Value = [double]"0.000786342660716947"
$ws.Range("P25").Value = [double]"0.0007863426607169469"
$ws.Range("Q25").Value = [double]"0.0003245010330000001"
$ws.Range("R25").Value = [double]"0.002920509297000001"
$ws.Range("S25").Value = [double]"1.237237639931013E-06"
$ws.Range("T25").Value = [double]"1.237237639931012E-06"
